# Segunda Atualização do Projeto
$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# ----------------------------------------------------------------------
# 1) " com a segurança " + bookmark(_GoBack) + "de todos" -> merge into a
#    single run of text " com a segurança de todos" and drop the bookmark
#    from this location (it is re-created later, near the map paragraph).
# ----------------------------------------------------------------------
$d.Content.Find.Execute("com a segurança de todos", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "com a segurança de todos", 2) | Out-Null

# remove the now-orphaned _GoBack bookmark (Word will recreate it below)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ----------------------------------------------------------------------
# 2) Remove <w:spacing w:val="-5"/> from the paragraph-mark run properties
#    of the paragraph that begins "Quando o usuário entrar no app...".
# ----------------------------------------------------------------------
$par = $d.Content.Find
$r = $d.Content
$r.Find.Execute("Quando o usuário entrar no app", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $r.Paragraphs(1).Range.Font.Spacing = 0
}

# ----------------------------------------------------------------------
# 3) Append new sentence after "...cadastrada por outros usuários." and
#    move the _GoBack bookmark to the end of the newly inserted text.
# ----------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("já cadastrada por outros usuários.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "já cadastrada por outros usuários ou pode adicionar manualmente a localização", 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("manualmente a localização", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
if ($r3.Find.Found) {
    $endRange = $r3.Duplicate
    $endRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $endRange)
}

# ----------------------------------------------------------------------
# 4) "Apresentar em um mapa relatos de outros usuários"
#    -> "Apresentar em um mapa com relatos de outros usuários"
# ----------------------------------------------------------------------
$d.Content.Find.Execute("Apresentar em um mapa relatos de outros usuários", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Apresentar em um mapa com relatos de outros usuários", 2) | Out-Null

# ----------------------------------------------------------------------
# 5) "...vários tipos de retos da população..." -> "...relatos da população..."
# ----------------------------------------------------------------------
$d.Content.Find.Execute("vários tipos de retos da população", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "vários tipos de relatos da população", 2) | Out-Null

# ----------------------------------------------------------------------
# 6) "...ao usuário desviar daquela localidade." -> "...ao usuário de evitar aquela localidade."
# ----------------------------------------------------------------------
$d.Content.Find.Execute("ao usuário desviar daquela localidade.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "ao usuário de evitar aquela localidade.", 2) | Out-Null

# ----------------------------------------------------------------------
# 7) "...também poderão conversar entre si através de um chat que o app
#    disponibilizará." -> "...também poderão interagir entre si através
#    do campo comentários."
# ----------------------------------------------------------------------
$d.Content.Find.Execute("também poderão conversar entre si através de um chat que o app disponibilizará.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "também poderão interagir entre si através do campo comentários.", 2) | Out-Null

# ----------------------------------------------------------------------
# 8) "...de forma online com off-line..." -> "...de forma online como off-line..."
# ----------------------------------------------------------------------
$d.Content.Find.Execute("de forma online com off-line", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "de forma online como off-line", 2) | Out-Null

# ----------------------------------------------------------------------
# 9) "...atualizará todos os outros usuários." -> "...atualizará para os outros usuários."
# ----------------------------------------------------------------------
$d.Content.Find.Execute("atualizará todos os outros usuários.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "atualizará para os outros usuários.", 2) | Out-Null

# ----------------------------------------------------------------------
# 10) "O usuário prestar seu" -> "O usuário pode prestar seu"
# ----------------------------------------------------------------------
$d.Content.Find.Execute("O usuário prestar seu", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "O usuário pode prestar seu", 2) | Out-Null

# ----------------------------------------------------------------------
# 11) "...torna-se fácil resolver determinados crimes..." -> "...torna-se fácil mostrar determinados crimes..."
# ----------------------------------------------------------------------
$d.Content.Find.Execute("torna-se fácil resolver determinados crimes", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "torna-se fácil mostrar determinados crimes", 2) | Out-Null

Write-Host "Edits applied"
